# Edit script: insert two new weekly price rows for "Perejil" (Vega Monumental
# Concepción) at the top of the data block (rows 82-83), shifting all the
# existing historical rows down by two (rows 82-181 -> 84-183), exactly as a
# new week of data is added on top of the existing logica_diaria series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 82; Excel automatically shifts every
# row from 82 downward to 84 onward (keeping all their existing data/format),
# and copies the formatting (e.g. the date style on column D) from the row
# directly above into the freshly inserted rows.
$ws.Rows("82:83").Insert()

# --- New row 82 ("Primera" quality) ---
$ws.Range("A82").Value = 11
$ws.Range("B82").Value = "Vega Monumental Concepción"
$ws.Range("C82").Value = "Bíobío"
$ws.Range("D82").Value = 44895
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = 100112044
$ws.Range("G82").Value = "Perejil"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 200
$ws.Range("K82").Value = 700
$ws.Range("L82").Value = 800
$ws.Range("M82").Value = 750
$ws.Range("N82").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O82").Value = "Región de Ñuble"
$ws.Range("P82").Value = 750
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"

# --- New row 83 ("Segunda" quality) ---
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44895
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100112044
$ws.Range("G83").Value = "Perejil"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Segunda"
$ws.Range("J83").Value = 100
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 600
$ws.Range("M83").Value = 600
$ws.Range("N83").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O83").Value = "Región de Ñuble"
$ws.Range("P83").Value = 600
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"

# Make sure the date cells keep/get the same date-time number format used by
# the rest of column D.
$ws.Range("D82:D83").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
